$d = $word.ActiveDocument

# --- Header date ---
$d.Content.Find.Execute("2025-06-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-23 Monday", 2) | Out-Null

# --- Table cell replacements (unique old values, safe to Find/Replace across the whole doc) ---
$replacements = @(
    @("77÷6=12, 5", "96÷9=10, 6"),
    @("50÷7=7, 1", "59÷4=14, 3"),
    @("11÷2=5, 1", "72÷7=10, 2"),
    @("94÷2=47, 0", "86÷8=10, 6"),
    @("54÷4=13, 2", "91÷4=22, 3"),
    @("49÷6=8, 1", "98÷2=49, 0"),
    @("14÷9=1, 5", "21÷2=10, 1"),
    @("42÷8=5, 2", "13÷8=1, 5"),
    @("41÷5=8, 1", "55÷9=6, 1"),
    @("99÷5=19, 4", "36÷8=4, 4"),
    @("93÷5=18, 3", "87÷6=14, 3"),
    @("94÷8=11, 6", "23÷8=2, 7"),
    @("40÷5=8, 0", "78÷4=19, 2"),
    @("79÷4=19, 3", "50÷6=8, 2"),
    @("96÷2=48, 0", "80÷2=40, 0"),
    @("33÷4=8, 1", "82÷3=27, 1"),
    @("17÷7=2, 3", "14÷5=2, 4"),
    @("75÷9=8, 3", "76÷9=8, 4"),
    @("64÷8=8, 0", "90÷8=11, 2"),
    @("21÷5=4, 1", "28÷7=4, 0"),
    @("71÷7=10, 1", "76÷3=25, 1"),
    @("32÷2=16, 0", "83÷8=10, 3"),
    @("51÷9=5, 6", "70÷8=8, 6")
)

foreach ($pair in $replacements) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null
}

# --- Duplicate value "13÷3=4, 1" appears twice in the table; address each cell directly ---
$table = $d.Tables.Item(1)
$table.Cell(9, 5).Range.Text = "72÷6=12, 0"
$table.Cell(13, 3).Range.Text = "13÷6=2, 1"
